$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finish convert relative coords to absolute coords:
# Row 3 (id=2): trigger type "arrow" -> "shoot"
$ws.Range("B3").Value = "shoot"
# Row 4 (id=3): trigger type "shoot" -> "penetrate"
$ws.Range("B4").Value = "penetrate"

# Update the active selection to B5
$ws.Range("B5").Select()
